$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(30).Insert()

$ws.Cells.Item(30, 1).Value = 11
$ws.Cells.Item(30, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(30, 3).Value = "Bíobío"
$ws.Cells.Item(30, 4).Value = 44972
$ws.Cells.Item(30, 5).Value = 8
$ws.Cells.Item(30, 6).Value = "Fruta"
$ws.Cells.Item(30, 7).Value = 100108
$ws.Cells.Item(30, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(30, 9).Value = 100108002
$ws.Cells.Item(30, 10).Value = "Mango"
$ws.Cells.Item(30, 11).Value = "Sin especificar"
$ws.Cells.Item(30, 12).Value = "Primera"
$ws.Cells.Item(30, 13).Value = 200
$ws.Cells.Item(30, 14).Value = 7500
$ws.Cells.Item(30, 15).Value = 8000
$ws.Cells.Item(30, 16).Value = 7750
$ws.Cells.Item(30, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(30, 18).Value = "Perú"
$ws.Cells.Item(30, 19).Value = 1938
$ws.Cells.Item(30, 20).Value = 4
